# Natmi following Dr Hou advice: refreshed Gdf11-Acvr2a ligand/receptor
# expression + specificity metrics (updated expressing-cell counts and
# derived edge-weight statistics) for rows 2-10.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = 3
$ws.Range("G2").Value = 1.296436
$ws.Range("H2").Value = 3.889308000000001
$ws.Range("I2").Value = 0.1944674516316147
$ws.Range("J2").Value = 0.1944674516316147
$ws.Range("K2").Value = 3
$ws.Range("M2").Value = 14.61878266666667
$ws.Range("N2").Value = 43.856348
$ws.Range("O2").Value = 0.2662829816142094
$ws.Range("P2").Value = 0.2662829816142094
$ws.Range("Q2").Value = 18.95231612524267
$ws.Range("R2").Value = 170.570845127184
$ws.Range("S2").Value = 0.05178337284738344
$ws.Range("T2").Value = 0.05178337284738344

$ws.Range("E3").Value = 3
$ws.Range("G3").Value = 1.296436
$ws.Range("H3").Value = 3.889308000000001
$ws.Range("I3").Value = 0.1944674516316147
$ws.Range("J3").Value = 0.1944674516316147
$ws.Range("K3").Value = 3
$ws.Range("M3").Value = 27.084169
$ws.Range("N3").Value = 81.25250700000001
$ws.Range("O3").Value = 0.4933415757187404
$ws.Range("P3").Value = 0.4933415757187404
$ws.Range("Q3").Value = 35.11289172168401
$ws.Range("R3").Value = 316.0160254951561
$ws.Range("S3").Value = 0.09593887901394875
$ws.Range("T3").Value = 0.09593887901394875

$ws.Range("E4").Value = 3
$ws.Range("G4").Value = 1.296436
$ws.Range("H4").Value = 3.889308000000001
$ws.Range("I4").Value = 0.1944674516316147
$ws.Range("J4").Value = 0.1944674516316147
$ws.Range("K4").Value = 3
$ws.Range("M4").Value = 13.19647366666667
$ws.Range("N4").Value = 39.589421
$ws.Range("O4").Value = 0.2403754426670501
$ws.Range("P4").Value = 0.2403754426670501
$ws.Range("Q4").Value = 17.10838353451867
$ws.Range("R4").Value = 153.975451810668
$ws.Range("S4").Value = 0.04674519977028255
$ws.Range("T4").Value = 0.04674519977028255

$ws.Range("E5").Value = 3
$ws.Range("G5").Value = 1.837275666666667
$ws.Range("H5").Value = 5.511827
$ws.Range("I5").Value = 0.2755942575194169
$ws.Range("J5").Value = 0.2755942575194169
$ws.Range("K5").Value = 3
$ws.Range("M5").Value = 14.61878266666667
$ws.Range("N5").Value = 43.856348
$ws.Range("O5").Value = 0.2662829816142094
$ws.Range("P5").Value = 0.2662829816142094
$ws.Range("Q5").Value = 26.85873366975511
$ws.Range("R5").Value = 241.728603027796
$ws.Range("S5").Value = 0.07338606060802458
$ws.Range("T5").Value = 0.07338606060802458

$ws.Range("E6").Value = 3
$ws.Range("G6").Value = 1.837275666666667
$ws.Range("H6").Value = 5.511827
$ws.Range("I6").Value = 0.2755942575194169
$ws.Range("J6").Value = 0.2755942575194169
$ws.Range("K6").Value = 3
$ws.Range("M6").Value = 27.084169
$ws.Range("N6").Value = 81.25250700000001
$ws.Range("O6").Value = 0.4933415757187404
$ws.Range("P6").Value = 0.4933415757187404
$ws.Range("Q6").Value = 49.76108465558767
$ws.Range("R6").Value = 447.8497619002891
$ws.Range("S6").Value = 0.1359621052636654
$ws.Range("T6").Value = 0.1359621052636654

$ws.Range("E7").Value = 3
$ws.Range("G7").Value = 1.837275666666667
$ws.Range("H7").Value = 5.511827
$ws.Range("I7").Value = 0.2755942575194169
$ws.Range("J7").Value = 0.2755942575194169
$ws.Range("K7").Value = 3
$ws.Range("M7").Value = 13.19647366666667
$ws.Range("N7").Value = 39.589421
$ws.Range("O7").Value = 0.2403754426670501
$ws.Range("P7").Value = 0.2403754426670501
$ws.Range("Q7").Value = 24.24555995357412
$ws.Range("R7").Value = 218.210039582167
$ws.Range("S7").Value = 0.06624609164772684
$ws.Range("T7").Value = 0.06624609164772682

$ws.Range("E8").Value = 3
$ws.Range("G8").Value = 3.532884666666666
$ws.Range("H8").Value = 10.598654
$ws.Range("I8").Value = 0.5299382908489685
$ws.Range("J8").Value = 0.5299382908489685
$ws.Range("K8").Value = 3
$ws.Range("M8").Value = 14.61878266666667
$ws.Range("N8").Value = 43.856348
$ws.Range("O8").Value = 0.2662829816142094
$ws.Range("P8").Value = 0.2662829816142094
$ws.Range("Q8").Value = 51.64647312839911
$ws.Range("R8").Value = 464.8182581555919
$ws.Range("S8").Value = 0.1411135481588014
$ws.Range("T8").Value = 0.1411135481588014

$ws.Range("E9").Value = 3
$ws.Range("G9").Value = 3.532884666666666
$ws.Range("H9").Value = 10.598654
$ws.Range("I9").Value = 0.5299382908489685
$ws.Range("J9").Value = 0.5299382908489685
$ws.Range("K9").Value = 3
$ws.Range("M9").Value = 27.084169
$ws.Range("N9").Value = 81.25250700000001
$ws.Range("O9").Value = 0.4933415757187404
$ws.Range("P9").Value = 0.4933415757187404
$ws.Range("Q9").Value = 95.68524536950866
$ws.Range("R9").Value = 861.167208325578
$ws.Range("S9").Value = 0.2614405914411262
$ws.Range("T9").Value = 0.2614405914411262

$ws.Range("E10").Value = 3
$ws.Range("G10").Value = 3.532884666666666
$ws.Range("H10").Value = 10.598654
$ws.Range("I10").Value = 0.5299382908489685
$ws.Range("J10").Value = 0.5299382908489685
$ws.Range("K10").Value = 3
$ws.Range("M10").Value = 13.19647366666667
$ws.Range("N10").Value = 39.589421
$ws.Range("O10").Value = 0.2403754426670501
$ws.Range("P10").Value = 0.2403754426670501
$ws.Range("Q10").Value = 46.62161947103711
$ws.Range("R10").Value = 419.594575239334
$ws.Range("S10").Value = 0.1273841512490408
$ws.Range("T10").Value = 0.1273841512490407
